$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 114
$ws1.Range("F5").Value = 6519
$ws1.Range("F6").Value = 78
$ws1.Range("F9").Value = 5883
$ws1.Range("F18").Value = 335
$ws1.Range("F21").Value = 4202
$ws1.Range("F22").Value = 34

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 114
$ws4.Range("F5").Value = 6519
$ws4.Range("F6").Value = 78
$ws4.Range("F9").Value = 5883
$ws4.Range("F18").Value = 335
$ws4.Range("F21").Value = 4202
$ws4.Range("F23").Value = 34
